$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Rename header G1 from "Practice Questions" to "Init Practice Questions"
$ws1.Range("G1").Value = "Init Practice Questions"

# Arrays row (row 2): Review Practice Questions count 0 -> 1
$ws1.Range("H2").Value = 1

# Graphs row (row 7): Own Implemention status Todo -> Done
$ws1.Range("F7").Value = "Done"

# Update selection to reflect the last active cell used (F6)
$ws1.Range("F6").Select()
